$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 0.1050794602663998
$ws.Range("B3").Value = 0.002616647741689461
$ws.Range("C3").Value = 0.000699325292216118
$ws.Range("D3").Value = 3.357613022016416
$ws.Range("E3").Value = 0.02823042075705754
$ws.Range("F3").Value = 0.0012459909502047
$ws.Range("G3").Value = 0.003987304533174221
$ws.Range("H3").Value = 0.1076961080080893
$ws.Range("B4").Value = 0.007992508331426796
$ws.Range("C4").Value = 0.001343919359700321
$ws.Range("D4").Value = 5.587098511468807
$ws.Range("E4").Value = 0.15382048003193
$ws.Range("F4").Value = 0.005358466861578125
$ws.Range("G4").Value = 0.01062654980127547
$ws.Range("H4").Value = 0.1130719685978266
$ws.Range("B5").Value = 0.01533011393984426
$ws.Range("C5").Value = 0.004068632681849238
$ws.Range("D5").Value = 4.517458142133467
$ws.Range("E5").Value = 0.08825169508336854
$ws.Range("F5").Value = 0.007355715875370269
$ws.Range("G5").Value = 0.02330451200431826
$ws.Range("H5").Value = 0.1204095742062441
$ws.Range("B6").Value = 0.02824243719058972
$ws.Range("C6").Value = 0.007633326614986345
$ws.Range("D6").Value = 5.726291775576585
$ws.Range("E6").Value = 0.06090806772698037
$ws.Range("F6").Value = 0.01328134870897334
$ws.Range("G6").Value = 0.0432035256722061
$ws.Range("H6").Value = 0.1333218974569895
$ws.Range("B7").Value = 0.02478967440452952
$ws.Range("C7").Value = 0.002817234013060411
$ws.Range("D7").Value = 6.108794103616588
$ws.Range("E7").Value = 0.05171976168954228
$ws.Range("F7").Value = 0.01926798021718982
$ws.Range("G7").Value = 0.03031136859186923
$ws.Range("H7").Value = 0.1298691346709293
$ws.Range("B8").Value = 0.02094605722373222
$ws.Range("C8").Value = 0.008042606108831698
$ws.Range("D8").Value = 6.123953916991939
$ws.Range("E8").Value = 0.2044865667363685
$ws.Range("F8").Value = 0.005182792912070262
$ws.Range("G8").Value = 0.03670932153539419
$ws.Range("H8").Value = 0.126025517490132
$ws.Range("B9").Value = 0.02085494399640334
$ws.Range("C9").Value = 0.00460145154192373
$ws.Range("D9").Value = 5.533274735784835
$ws.Range("E9").Value = 0.1057728683738496
$ws.Range("F9").Value = 0.01183623576427352
$ws.Range("G9").Value = 0.02987365222853316
$ws.Range("H9").Value = 0.1259344042628031
$ws.Range("B10").Value = -0.1050794602663998
$ws.Range("C10").Value = 0.0005069664040439407
$ws.Range("D10").Value = -228.5873401718663
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = -0.106073099310044
$ws.Range("G10").Value = -0.1040858212227555
$ws.Range("B11").Value = -0.04673655610730883
$ws.Range("C11").Value = 0.0005441861370308571
$ws.Range("D11").Value = -93.21436166117368
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = -0.04780314471529819
$ws.Range("G11").Value = -0.04566996749931947
$ws.Range("H11").Value = 0.05834290415909097
$ws.Range("B12").Value = -0.03257552022828045
$ws.Range("C12").Value = 0.0005257913416741138
$ws.Range("D12").Value = -68.43658122450861
$ws.Range("E12").Value = [double]"6.969491333696976e-74"
$ws.Range("F12").Value = -0.03360605559580012
$ws.Range("G12").Value = -0.03154498486076079
$ws.Range("H12").Value = 0.07250394003811936
$ws.Range("B13").Value = -0.03046910075102528
$ws.Range("C13").Value = 0.000518041775051201
$ws.Range("D13").Value = -65.21076060503125
$ws.Range("E13").Value = [double]"2.416895154381112e-60"
$ws.Range("F13").Value = -0.03148444720006648
$ws.Range("G13").Value = -0.02945375430198409
$ws.Range("H13").Value = 0.07461035951537452
$ws.Range("B14").Value = -0.02381916608994797
$ws.Range("C14").Value = 0.0005118445595539397
$ws.Range("D14").Value = -53.34972717351289
$ws.Range("E14").Value = [double]"6.660663760506254e-06"
$ws.Range("F14").Value = -0.02482236617903563
$ws.Range("G14").Value = -0.02281596600086032
$ws.Range("H14").Value = 0.08126029417645184
$ws.Range("B15").Value = -0.02212653754626784
$ws.Range("C15").Value = 0.0005041223902920719
$ws.Range("D15").Value = -50.02791680169585
$ws.Range("E15").Value = [double]"3.924495661535669e-19"
$ws.Range("F15").Value = -0.02311460241235599
$ws.Range("G15").Value = -0.02113847268017969
$ws.Range("H15").Value = 0.08295292272013197
$ws.Range("B16").Value = -0.0213452221449854
$ws.Range("C16").Value = 0.0005019659167107206
$ws.Range("D16").Value = -48.69624393507612
$ws.Range("E16").Value = [double]"3.849102762630352e-23"
$ws.Range("F16").Value = -0.02232906038915536
$ws.Range("G16").Value = -0.02036138390081545
$ws.Range("H16").Value = 0.08373423812141439
$ws.Range("B17").Value = -0.02079176833113184
$ws.Range("C17").Value = 0.0005074026780684215
$ws.Range("D17").Value = -47.26461036401172
$ws.Range("E17").Value = 0.001786222043244245
$ws.Range("F17").Value = -0.02178626247006482
$ws.Range("G17").Value = -0.01979727419219885
$ws.Range("H17").Value = 0.08428769193526797
$ws.Range("B18").Value = -0.01916729019050374
$ws.Range("C18").Value = 0.000509892705822016
$ws.Range("D18").Value = -42.37822927458218
$ws.Range("E18").Value = [double]"4.244517655294516e-07"
$ws.Range("F18").Value = -0.02016666470781468
$ws.Range("G18").Value = -0.01816791567319281
$ws.Range("H18").Value = 0.08591217007589606
$ws.Range("B19").Value = -0.01423106186205974
$ws.Range("C19").Value = 0.0005047539611743892
$ws.Range("D19").Value = -31.38864794155348
$ws.Range("E19").Value = 0.04931255358225387
$ws.Range("F19").Value = -0.01522036459181895
$ws.Range("G19").Value = -0.01324175913230053
$ws.Range("H19").Value = 0.09084839840434007
$ws.Range("B20").Value = -0.01120274946867032
$ws.Range("C20").Value = 0.0005153325585861909
$ws.Range("D20").Value = -24.3023365592085
$ws.Range("E20").Value = 0.06494468931339363
$ws.Range("F20").Value = -0.01221278592917946
$ws.Range("G20").Value = -0.01019271300816118
$ws.Range("H20").Value = 0.09387671079772948
$ws.Range("B21").Value = -0.00889528363895584
$ws.Range("C21").Value = 0.0005188113001432928
$ws.Range("D21").Value = -18.28262899442666
$ws.Range("E21").Value = 0.01281937232005859
$ws.Range("F21").Value = -0.009912138328800366
$ws.Range("G21").Value = -0.007878428949111315
$ws.Range("H21").Value = 0.09618417662744397
$ws.Range("B22").Value = -0.006696284553999147
$ws.Range("C22").Value = 0.0005183352524455296
$ws.Range("D22").Value = -13.18383853313903
$ws.Range("E22").Value = 0.1150947285512146
$ws.Range("F22").Value = -0.007712206206320567
$ws.Range("G22").Value = -0.005680362901677728
$ws.Range("H22").Value = 0.09838317571240066
$ws.Range("B23").Value = -0.006085875672799853
$ws.Range("C23").Value = 0.0005197762041982257
$ws.Range("D23").Value = -12.52915930639699
$ws.Range("E23").Value = 0.0660059817044503
$ws.Range("F23").Value = -0.007104621554486949
$ws.Range("G23").Value = -0.005067129791112759
$ws.Range("H23").Value = 0.09899358459359996
$ws.Range("B24").Value = -0.003915038677801314
$ws.Range("C24").Value = 0.0005208360296117919
$ws.Range("D24").Value = -7.519050021944413
$ws.Range("E24").Value = 0.04840239168073581
$ws.Range("F24").Value = -0.004935861773897806
$ws.Range("G24").Value = -0.002894215581704823
$ws.Range("H24").Value = 0.1011644215885985
$ws.Range("B25").Value = -0.002740078882508847
$ws.Range("C25").Value = 0.0005129720100060278
$ws.Range("D25").Value = -6.064298565175729
$ws.Range("E25").Value = 0.1219549883686271
$ws.Range("F25").Value = -0.003745488732031723
$ws.Range("G25").Value = -0.001734669032985971
$ws.Range("H25").Value = 0.102339381383891
$ws.Range("B26").Value = 0.02383908531155556
$ws.Range("C26").Value = 0.00455665840254458
$ws.Range("D26").Value = 13.5958273921253
$ws.Range("E26").Value = 0.05433002101649512
$ws.Range("F26").Value = 0.01328134870897334
$ws.Range("G26").Value = 0.03276999761994104
$ws.Range("H26").Value = 0.1289185455779553
